# Add a new worksheet "OHT Relation" after the last existing sheet and
# populate it with the OHT id relation matrix, then make it the active sheet.

$wb = $excel.ActiveWorkbook

$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)

$ohtSheet = $wb.Worksheets.Add($null, $lastSheet)
$ohtSheet.Name = "OHT Relation"

# Header row: top-left label + column ids 0,1,2
$ohtSheet.Range("A1").Value = "OHT_id"
$ohtSheet.Range("B1").Value = 0
$ohtSheet.Range("C1").Value = 1
$ohtSheet.Range("D1").Value = 2

# Row id 0
$ohtSheet.Range("A2").Value = 0
$ohtSheet.Range("B2").Value = 0
$ohtSheet.Range("C2").Value = 1
$ohtSheet.Range("D2").Value = 0

# Row id 1
$ohtSheet.Range("A3").Value = 1
$ohtSheet.Range("B3").Value = -1
$ohtSheet.Range("C3").Value = 0
$ohtSheet.Range("D3").Value = 0

# Row id 2
$ohtSheet.Range("A4").Value = 2
$ohtSheet.Range("B4").Value = 0
$ohtSheet.Range("C4").Value = 0
$ohtSheet.Range("D4").Value = 0

# Page setup to match a plain portrait A4 sheet
$pageSetup = $ohtSheet.PageSetup
$pageSetup.PaperSize = 9
$pageSetup.Orientation = 1

# Leave the final selection on D3, and make this new sheet the active tab
$null = $ohtSheet.Range("D3").Select()
$null = $ohtSheet.Activate()
